$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title / author / email -------------------------------------------------
Replace-Text "Microcosm of Life: Exploring the Wonders of a Single Cell" "The Marvelous World of Biology: Unraveling the Secrets of Life"
Replace-Text "Amelia Lee" "Amelia Becker"
Replace-Text "lee@biosphere" "becker@validmail"

# --- First body paragraph ----------------------------------------------------
Replace-Text "In the vast tapestry of life, from the towering sequoias to the microscopic plankton, there exists an intricate world within each living cell" "Immerse yourself in the enchanting realm of biology, an intricate science that delves into the very essence of life"

Replace-Text " The study of cytology delves into this cellular microcosm, revealing a symphony of complex processes, intricate structures, and remarkable adaptations that collectively give rise to the diversity of life on Earth" " Prepare to embark on a captivating journey as we explore the microscopic universe that governs all living organisms"

Replace-Text " As we embark on this journey into the depths of a single cell, we discover a universe teeming with wonder, innovation, and unity" " Biology unveils the intricate symphony of life, from the graceful dance of molecules to the intricate web of ecosystems"

# New sentence inserted after the previous one (and after the existing "." run)
$rng = $d.Content
$rng.Find.Execute(" Biology unveils the intricate symphony of life, from the graceful dance of molecules to the intricate web of ecosystems") | Out-Null
$insPoint = $d.Range($rng.End, $rng.End)
$insPoint.Collapse(0)
$insPoint.MoveStart(1, 1) | Out-Null
$insPoint.InsertAfter(" Join us as we unveil the secrets hidden within the tapestry of life, deciphering the enigmatic code that governs our existence")

Replace-Text "Like a miniature city, a cell is a bustling hub of activity" "Biology unveils the intricate interplay between organisms and their environment, revealing the delicate balance that sustains life"

Replace-Text " Organelles, each with specialized functions, perform their intricate tasks, ensuring the cell's survival and growth" " Witness the awe-inspiring diversity of species, each possessing unique adaptations that enable them to thrive in an ever-changing world"

Replace-Text " The nucleus, the control center, houses the genetic blueprint that orchestrates the cell's identity and function" " Explore the fascinating world of genetics, uncovering the blueprint of life and comprehending the mechanisms that govern inheritance and variation"

Replace-Text " The mitochondria, the energy powerhouses, generate the fuel that drives cellular processes" " Immerse yourself in the dynamic world of ecosystems, where organisms interact in intricate ways, shaping and sustaining the intricate web of life"

# Remove the three sentences that followed (endoplasmic reticulum / lysosomes), keeping the trailing "."
Replace-Text ". The endoplasmic reticulum, a vast network of membranes, facilitates the transport and synthesis of molecules. The lysosomes, the cellular recycling centers, break down waste and cellular debris" ""

Replace-Text "Cells communicate with each other through a sophisticated signaling network, exchanging vital information that coordinates their actions and maintains tissue homeostasis" "Biology ignites a sense of wonder and curiosity, inspiring us to explore the complexities of the natural world"

Replace-Text " They respond to external stimuli, such as changes in temperature, pH, or nutrient availability, adjusting their behavior to adapt to a dynamic environment" " Whether examining the delicate structure of a cell or delving into the vastness of an ecosystem, biology offers a profound understanding of the world we inhabit"

Replace-Text " In multicellular organisms, cells work together in harmony, forming tissues, organs, and systems that carry out specialized functions, demonstrating the remarkable power of collective action" " It challenges us to question, investigate, and unravel the mysteries that surround us, leading us on a lifelong journey of discovery and enlightenment"

# --- Summary paragraph --------------------------------------------------------
Replace-Text "The exploration of a single cell unveils a microcosm of life, showcasing the extraordinary complexity, diversity, and unity that characterize all living organisms" "Biology unveils the intricate tapestry of life, revealing the breathtaking diversity of organisms and the delicate balance of ecosystems"

Replace-Text " Within the confines of this tiny world, we witness the intricate interplay of organelles, the marvels of genetic regulation, and the remarkable adaptability of cellular life" " It explores the inner workings of cells, the secrets of genetics, and the intricate mechanisms that govern inheritance and variation"

Replace-Text " The study of cytology provides a window into the fundamental processes that govern life, offering insights into our own biology, health, and the interconnectedness of all living things" " Biology ignites a sense of wonder and curiosity, inspiring us to explore the complexities of the natural world and unravel the mysteries that surround us"

# --- New trailing empty paragraph --------------------------------------------
$d.Content.InsertParagraphAfter()
